$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old row 5 (26KBF1KVM004 / radha) -- rows below shift up
$ws.Rows("5").Delete()

# --- New card codes, column A, rows 6-19 (entered first as a block) ---
$ws.Range("A6").Value  = "26KBF1KVM005"
$ws.Range("A7").Value  = "26KBF1KVM006"
$ws.Range("A8").Value  = "26KBF1KVM007"
$ws.Range("A9").Value  = "26KBF1KVM008"
$ws.Range("A10").Value = "26KBF1KVM009"
$ws.Range("A11").Value = "26KBF1KVM010"
$ws.Range("A12").Value = "26KBF1KVM011"
$ws.Range("A13").Value = "26KBF1KVM012"
$ws.Range("A14").Value = "26KBF1KVM013"
$ws.Range("A15").Value = "26KBF1KVM014"
$ws.Range("A16").Value = "26KBF1KVM015"
$ws.Range("A17").Value = "26KBF1KVM016"
$ws.Range("A18").Value = "26KBF1KVM017"
$ws.Range("A19").Value = "26KBF1KVM018"

# Amount paid for the first new row
$ws.Range("C6").Value = 2000

# --- Names, column B, rows 6-15, then 17-18, then 16, then 19 ---
$ws.Range("B6").Value  = "Maalamma"
$ws.Range("B7").Value  = "Lakshmi Mahadevamma 1"
$ws.Range("B8").Value  = "Manjanna 1"
$ws.Range("B9").Value  = "Manjanna 2"
$ws.Range("B10").Value = "Manjanna 3"
$ws.Range("B11").Value = "Shivanna"
$ws.Range("B12").Value = "Srinivas leelamma 1"
$ws.Range("B13").Value = "Srinivas leelamma 2"
$ws.Range("B14").Value = "Srinivas leelamma 3"
$ws.Range("B15").Value = "Srinivas leelamma 4"
$ws.Range("B17").Value = "Suma Mahadevamma 1"
$ws.Range("B18").Value = "Suma Mahadevamma 2"
$ws.Range("B16").Value = "Srinivas leelamma 5"
$ws.Range("B19").Value = "Suma Muddanayaka"

# --- New card codes, column A, rows 20-23 ---
$ws.Range("A20").Value = "26KBF1KVM019"
$ws.Range("A21").Value = "26KBF1KVM020"
$ws.Range("A22").Value = "26KBF1KVM021"
$ws.Range("A23").Value = "26KBF1KVM022"

# --- Names, column B, rows 20-21 ---
$ws.Range("B20").Value = "Muddanayaka Mahadevi"
$ws.Range("B21").Value = "Bhagya "

# Column B widened (content got longer) -- closest attainable step to the
# recorded 21.77734375 OOXML width.
$ws.Columns("B").ColumnWidth = 21

# Final selection lands on B21, matching the author's last edited cell
$ws.Range("B21").Select()
